# Update of the spreadsheet
# Applies the content changes described by the commit: several parameter
# labels in column A gain explicit unit suffixes, a new "v0.1" value is
# added at B41, the "Version of the spreadsheet" note is reworded, and the
# active selection moves to G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A parameter-name labels that now carry their unit in the name.
$ws.Range("A6").Value  = "Exposure_s"
$ws.Range("A9").Value  = "Sampling_Step_Size_microm(dx,dy,dz)"
$ws.Range("A10").Value = "Field_Of_View_microm(X,Y,Z)"
$ws.Range("A16").Value = "Confocal_Pinhole_Diameter_AU"
$ws.Range("A23").Value = "Illumination_Power_mW"
$ws.Range("A26").Value = "Laser_Drift_MHz_per_h"
$ws.Range("A29").Value = "Scan_Amplitude_GHz"
$ws.Range("A31").Value = "Scattering_Angle_deg"
$ws.Range("A32").Value = "Spectral_Resolution_MHz"
$ws.Range("A33").Value = "x-Mechanical_Resolution_microm"
$ws.Range("A34").Value = "x-Optical_Resolution_microm"
$ws.Range("A35").Value = "y-Mechanical_Resolution_microm"
$ws.Range("A36").Value = "y-Optical_Resolution_microm"
$ws.Range("A37").Value = "z-Mechanical_Resolution_microm"
$ws.Range("A38").Value = "z-Optical_Resolution_microm"

# New version value cell, and reworded note about it.
$ws.Range("B41").Value = "v0.1"
$ws.Range("D41").Value = "Version of the spreadsheet - Don't change"

# Move / update the active selection shown when the sheet is reopened.
$ws.Range("G13").Select()
